$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '41.527.52'
$ws.Cells.Item(2, 5).Value = '  -0.23%  '
$ws.Cells.Item(3, 4).Value = '2.194.74'
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '228.99'
$cell.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.86%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.615'
$cell.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -4.48%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '59.87'
$cell.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -6.92%  '
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 5).Value = '  -3.31%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '56.70'
$cell.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -5.79%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0882'
$cell.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -2.55%  '
$ws.Cells.Item(12, 5).Value = '  -2.03%  '
$ws.Cells.Item(13, 4).Value = '2.521.70'
$ws.Cells.Item(13, 5).Value = '  -2.54%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '15.28'
$cell.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -5.68%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '22.09'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.62'
$cell.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -1.30%  '
$ws.Cells.Item(17, 5).Value = '  -5.21%  '
$ws.Cells.Item(18, 4).Value = '2.204.34'
$ws.Cells.Item(18, 5).Value = '  -1.91%  '
$ws.Cells.Item(19, 4).Value = '41.461.52'
$ws.Cells.Item(19, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '71.70'
$cell.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -3.10%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0892'
$ws.Cells.Item(21, 5).Value = '  -4.12%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.00'
$cell.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -3.26%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '241.42'
$cell.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -4.55%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.16%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.34'
$cell.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -2.73%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.28'
$cell.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -2.50%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '9.54'
$cell.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -3.65%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '168.64'
$cell.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -2.15%  '
$ws.Cells.Item(29, 5).Value = '  -7.08%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.44'
$cell.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -1.11%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '19.62'
$cell.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -4.32%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.57'
$cell.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -8.82%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.120'
$cell.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -4.29%  '
$ws.Cells.Item(34, 5).Value = '  -2.77%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.59'
$cell.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -4.01%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0642'
$cell.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +0.13%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.29'
$cell.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -9.00%  '
$ws.Cells.Item(38, 5).Value = '  -5.16%  '
$ws.Cells.Item(39, 5).Value = '  -9.08%  '
$ws.Cells.Item(40, 2).Value = 'BinanceUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.25%  '
$ws.Cells.Item(41, 2).Value = 'TerraClassic'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.000238'
$cell.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -7.75%  '
$ws.Cells.Item(42, 5).Value = '  -2.93%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '8.44'
$cell.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -4.38%  '
$ws.Cells.Item(44, 5).Value = '  -6.03%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '96.44'
$cell.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -6.17%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.18'
$cell.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -4.43%  '
$ws.Cells.Item(47, 4).Value = '1.456.61'
$ws.Cells.Item(47, 5).Value = '  -3.68%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.32'
$cell.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -16.25%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '16.26'
$cell.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -8.52%  '
$ws.Cells.Item(50, 5).Value = '  -1.76%  '
$ws.Cells.Item(51, 5).Value = '  -7.33%  '
